# Insert a new row at row 280 (shifts existing rows 280:349 down to 281:350)
# and populate the newly inserted row with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A280:R280").Insert()

$ws.Range("A280").Value = 4
$ws.Range("B280").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C280").Value = "Los Lagos"
$ws.Range("D280").Value = 44964
$ws.Range("E280").Value = 10
$ws.Range("F280").Value = 100112032
$ws.Range("G280").Value = "Zapallo italiano"
$ws.Range("H280").Value = "Sin especificar"
$ws.Range("I280").Value = "Primera"
$ws.Range("J280").Value = 250
$ws.Range("K280").Value = 11000
$ws.Range("L280").Value = 11000
$ws.Range("M280").Value = 11000
$ws.Range("N280").Value = "$/caja 50 unidades"
$ws.Range("O280").Value = "Región de O'Higgins"
$ws.Range("P280").Value = 220
$ws.Range("Q280").Value = 50
$ws.Range("R280").Value = "Hortaliza"
